$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("A209").Value = 44776
Write-Output ($ws.Range("A209").Value2)
